# Update "想去人数" (interest count) figures for the affected events.
# These values live both on the "展览" sheet and on the aggregated
# "全部类型" sheet, which mirrors the same rows.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 1652   # was 1651
    $ws.Range("F5").Value = 1065   # was 1064
    $ws.Range("F6").Value = 641    # was 633
    $ws.Range("F8").Value = 5753   # was 5748
}
